# Updated cryptos list on Sun Oct  8 20:19:57 UTC 2023 with GitHub Actions
# Refresh the "Price" (D) and "Volume(1h)" (E) columns on the active sheet
# to match the latest coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.907.38"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.635.49"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.32"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("E9").Value = "  -0.30%  "
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0883"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "1.867.51"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").Value = "1.642.31"
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("E14").Value = "  -1.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.562"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").Value = "27.919.17"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.75%  "
$ws.Range("D20").Value = "0.0₃0719"
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("E23").Value = "  -2.30%  "
$ws.Range("E24").Value = "  +0.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.36%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0481"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.10"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.35%  "
$ws.Range("D34").Value = "1.399.54"
$ws.Range("E34").Value = "  +0.37%  "
$ws.Range("E35").Value = "  +3.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.94%  "
$ws.Range("E37").Value = "  -0.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0170"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.558"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.851"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.11%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("E42").Value = "  -1.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "66.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("D46").Value = "1.776.19"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.15%  "
$ws.Range("E49").Value = "  +2.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0504"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.58%  "
